# "testing results updated for week 7"
# - Week 7 results: MIN(a) @ NYG(h) game result is now in: NYG(h) won.
#   Header cell S1 updates to rich text with the winning team (NYG(h))
#   shown in bold red, and every model's "Wrong" tally (column C) goes
#   up by one since all the models had predicted MIN to win.
# - Season results becomes the active sheet/tab again (its formulas in
#   B:G recalc automatically from the updated Week 7 numbers).

$wb = $excel.ActiveWorkbook

$wk7 = $wb.Worksheets.Item("Week 7 results")
$wk7.Activate()

# Mark the final result for the MIN(a) @ NYG(h) game: NYG(h) won.
$s1 = $wk7.Range("S1")
$s1.Value = "MIN(a) @ NYG(h)"
$winner = $s1.Characters(10, 6)
$winner.Font.Bold = $true
$winner.Font.Color = 255

# Every model guessed this game wrong (they all favored MIN), so bump
# each row's "Wrong" column (C) by one.
$wk7.Range("C2").Value = 7
$wk7.Range("C3").Value = 5
$wk7.Range("C4").Value = 5
$wk7.Range("C5").Value = 5
$wk7.Range("C6").Value = 6
$wk7.Range("C7").Value = 6

# Leave the Week 7 sheet's own selection on C6 (matches the author's
# last click while entering these results).
$wk7.Range("C6").Select()

# Season results is the sheet that should be active/selected when the
# workbook is reopened.
$season = $wb.Worksheets.Item("Season results")
$season.Activate()
